$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Cant. Periodos" from 2 to 1
$ws.Range("F13").Value = 1

# Update "Valor Mora" total from 113880 to 56940
$ws.Range("E11").Value = 56940

# Remove the duplicate/old period row (period 2508), shifting everything below up
$ws.Rows("17:17").Delete()
